$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Maps" to "Map"
$ws.Name = "Map"

# Clear the existing data region before rewriting it
$ws.Range("A1:B7").ClearContents()

# Rebuild the table: Country/Cities header, Lietuva row with its cities
$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Cities"
$ws.Range("A2").Value = "Lietuva"
$ws.Range("B2").Value = "Vilnius"
$ws.Range("B3").Value = "Kaunas"
$ws.Range("B4").Value = "Siauliai"
$ws.Range("B5").Value = "Jonava"
$ws.Range("B6").Value = "Klaipeda"

# Move the active selection to A3, matching the saved view state
$ws.Range("A3").Select()
